# Fixed some star locations and magnitudes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("starData")

# --- Data corrections -------------------------------------------------

# u Ori (row 3): magnitude corrected
$ws.Range("G3").Value = 2.75

# Saif al Jabbar (row 9): RA/Dec had been copy-pasted from the row above;
# put in the real values.
$ws.Range("E9").Value = 0.22534722222222223
$ws.Range("F9").FormulaArray = "=-SUMPRODUCT(60^{0,-1,-2},{2,23,47})"

# --- Remove now-unused formatting on column D / column E --------------

$ws.Range("D1:D16").ClearFormats()
$ws.Columns("E").ClearFormats()

# Clearing the whole column leaves stray blank cells on the other rows
# that already had content further down the sheet; drop those again.
$ws.Range("E17").Clear()
$ws.Range("E26").Clear()
$ws.Range("E34").Clear()

# Re-apply the formatting that column E's cells still need (only the
# column's own bestFit/override is being dropped, not the individual
# cell styles).
$ws.Range("E1").HorizontalAlignment = -4152   # xlRight
$ws.Range("E2:E15").HorizontalAlignment = -4152   # xlRight
$ws.Range("E2:E15").NumberFormat = "h:mm:ss"

# --- Drop the stray blank row 16 ---------------------------------------

$ws.Range("D16:E16").Clear()

# --- Selection moved to G4 ---------------------------------------------

$ws.Range("G4").Select()
